$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (F column) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 6
$wsExpo.Range("F5").Value = 3811
$wsExpo.Range("F7").Value = 49
$wsExpo.Range("F8").Value = 229
$wsExpo.Range("F9").Value = 13

# Sheet "全部类型" - update 想去人数 (F column) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 6
$wsAll.Range("F9").Value = 3811
$wsAll.Range("F11").Value = 49
$wsAll.Range("F13").Value = 229
$wsAll.Range("F14").Value = 13
